{"js": "// Replace the date and each three-digit x one-digit multiplication\n// problem in the practice table with the new values from the target\n// revision. Every old value is unique in the document, so a simple\n// search-and-replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"2024-03-04 Monday\", \"2024-03-05 Tuesday\"],\n  [\"884\u00d76=5304\", \"475\u00d74=1900\"],\n  [\"471\u00d73=1413\", \"430\u00d72=860\"],\n  [\"862\u00d78=6896\", \"629\u00d77=4403\"],\n  [\"317\u00d73=951\", \"581\u00d75=2905\"],\n  [\"940\u00d77=6580\", \"221\u00d77=1547\"],\n  [\"900\u00d79=8100\", \"302\u00d72=604\"],\n  [\"492\u00d72=984\", \"649\u00d77=4543\"],\n  [\"847\u00d78=6776\", \"948\u00d78=7584\"],\n  [\"756\u00d73=2268\", \"361\u00d79=3249\"],\n  [\"351\u00d74=1404\", \"412\u00d79=3708\"],\n  [\"254\u00d72=508\", \"441\u00d76=2646\"],\n  [\"875\u00d74=3500\", \"747\u00d73=2241\"],\n  [\"565\u00d77=3955\", \"950\u00d78=7600\"],\n  [\"474\u00d72=948\", \"127\u00d79=1143\"],\n  [\"765\u00d72=1530\", \"973\u00d79=8757\"],\n  [\"832\u00d77=5824\", \"725\u00d78=5800\"],\n  [\"522\u00d72=1044\", \"666\u00d72=1332\"],\n  [\"860\u00d73=2580\", \"814\u00d72=1628\"],\n  [\"132\u00d74=528\", \"561\u00d74=2244\"],\n  [\"303\u00d76=1818\", \"648\u00d77=4536\"],\n  [\"426\u00d79=3834\", \"934\u00d77=6538\"],\n  [\"128\u00d74=512\", \"959\u00d73=2877\"],\n  [\"125\u00d79=1125\", \"274\u00d79=2466\"],\n  [\"901\u00d72=1802\", \"478\u00d73=1434\"],\n  [\"260\u00d76=1560\", \"811\u00d74=3244\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each three-digit x one-digit multiplication\n# problem in the practice table with the new values from the target\n# revision. Every old value is unique in the document, so a simple\n# Find/Replace per pair is safe and unambiguous.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n    @(\"2024-03-04 Monday\", \"2024-03-05 Tuesday\"),\n    @(\"884\u00d76=5304\", \"475\u00d74=1900\"),\n    @(\"471\u00d73=1413\", \"430\u00d72=860\"),\n    @(\"862\u00d78=6896\", \"629\u00d77=4403\"),\n    @(\"317\u00d73=951\", \"581\u00d75=2905\"),\n    @(\"940\u00d77=6580\", \"221\u00d77=1547\"),\n    @(\"900\u00d79=8100\", \"302\u00d72=604\"),\n    @(\"492\u00d72=984\", \"649\u00d77=4543\"),\n    @(\"847\u00d78=6776\", \"948\u00d78=7584\"),\n    @(\"756\u00d73=2268\", \"361\u00d79=3249\"),\n    @(\"351\u00d74=1404\", \"412\u00d79=3708\"),\n    @(\"254\u00d72=508\", \"441\u00d76=2646\"),\n    @(\"875\u00d74=3500\", \"747\u00d73=2241\"),\n    @(\"565\u00d77=3955\", \"950\u00d78=7600\"),\n    @(\"474\u00d72=948\", \"127\u00d79=1143\"),\n    @(\"765\u00d72=1530\", \"973\u00d79=8757\"),\n    @(\"832\u00d77=5824\", \"725\u00d78=5800\"),\n    @(\"522\u00d72=1044\", \"666\u00d72=1332\"),\n    @(\"860\u00d73=2580\", \"814\u00d72=1628\"),\n    @(\"132\u00d74=528\", \"561\u00d74=2244\"),\n    @(\"303\u00d76=1818\", \"648\u00d77=4536\"),\n    @(\"426\u00d79=3834\", \"934\u00d77=6538\"),\n    @(\"128\u00d74=512\", \"959\u00d73=2877\"),\n    @(\"125\u00d79=1125\", \"274\u00d79=2466\"),\n    @(\"901\u00d72=1802\", \"478\u00d73=1434\"),\n    @(\"260\u00d76=1560\", \"811\u00d74=3244\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n"}
